$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 173, pushing the existing
# rows 173:274 down to 174:275.
$ws.Rows("173:173").Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(173, 1).Value = 7
$ws.Cells.Item(173, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(173, 3).Value = "Ñuble"
$ws.Cells.Item(173, 4).Value = 44830
$ws.Cells.Item(173, 5).Value = 16
$ws.Cells.Item(173, 6).Value = 100112003
$ws.Cells.Item(173, 7).Value = "Ajo"
$ws.Cells.Item(173, 8).Value = "Chino"
$ws.Cells.Item(173, 9).Value = "Primera"
$ws.Cells.Item(173, 10).Value = 60
$ws.Cells.Item(173, 11).Value = 21000
$ws.Cells.Item(173, 12).Value = 22000
$ws.Cells.Item(173, 13).Value = 21500
$ws.Cells.Item(173, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(173, 15).Value = "China"
$ws.Cells.Item(173, 16).Value = 2150
$ws.Cells.Item(173, 17).Value = 10
$ws.Cells.Item(173, 18).Value = "Hortaliza"
